# Apply the updated cryptocurrency price/volume snapshot described in the commit diff.
# A helper is used for the "Price" (column D) cells: several of the new prices are plain
# decimal numbers (e.g. "0.9998", "113.40"). Every Price/Volume cell in this sheet is
# stored as *text*, so a leading apostrophe is used to stop Excel from re-interpreting
# these as numbers (which would also strip meaningful trailing zeros / change cell type).
function Set-TextValue($range, [string]$text) {
    $range.Value = "`'$text"
}

$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

# Row 2 (Bitcoin)
$ws.Range('D2').Value = '27.680.91'
$ws.Range('E2').Value = '  -0.15%  '

# Row 3 (Ethereum)
$ws.Range('D3').Value = '1.900.35'
$ws.Range('E3').Value = '  +0.47%  '

# Row 4 (TetherUSD)
Set-TextValue $ws.Range('D4') '0.9998'
$ws.Range('E4').Value = '  -0.60%  '

# Row 5 (BNB)
Set-TextValue $ws.Range('D5') '312.01'
$ws.Range('E5').Value = '  -0.85%  '

# Row 6 (USDC)
Set-TextValue $ws.Range('D6') '0.9998'
$ws.Range('E6').Value = '  -0.54%  '

# Row 7 (XRP)
Set-TextValue $ws.Range('D7') '0.5184'
$ws.Range('E7').Value = '  +8.04%  '

# Row 8 (Cardano)
Set-TextValue $ws.Range('D8') '0.3779'
$ws.Range('E8').Value = '  -0.38%  '

# Row 9 (Dogecoin)
Set-TextValue $ws.Range('D9') '0.07237'
$ws.Range('E9').Value = '  -1.19%  '

# Row 10 (Solana)
Set-TextValue $ws.Range('D10') '20.98'
$ws.Range('E10').Value = '  +2.66%  '

# Row 11 (Polygon)
Set-TextValue $ws.Range('D11') '0.8938'
$ws.Range('E11').Value = '  -2.71%  '

# Row 12 (TRON)
$ws.Range('E12').Value = '  -0.76%  '

# Row 13 (WrappedEther)
$ws.Range('D13').Value = '1.899.06'
$ws.Range('E13').Value = '  +0.16%  '

# Row 14 (Polkadot)
Set-TextValue $ws.Range('D14') '5.439'
$ws.Range('E14').Value = '  -0.34%  '

# Row 15 (Litecoin)
Set-TextValue $ws.Range('D15') '92.03'
$ws.Range('E15').Value = '  +1.10%  '

# Row 16 (BinanceUSD)
Set-TextValue $ws.Range('D16') '1.000'
$ws.Range('E16').Value = '  -0.57%  '

# Row 17 (ShibaInu)
Set-TextValue $ws.Range('D17') '0.000008713'
$ws.Range('E17').Value = '  -0.79%  '

# Row 18 (Dai)
$ws.Range('E18').Value = '  -0.25%  '

# Row 19 (WrappedBTC)
$ws.Range('D19').Value = '27.715.20'
$ws.Range('E19').Value = '  -0.26%  '

# Row 20 (Avalanche)
Set-TextValue $ws.Range('D20') '14.45'
$ws.Range('E20').Value = '  -0.37%  '

# Row 21 (Uniswap)
Set-TextValue $ws.Range('D21') '5.137'
$ws.Range('E21').Value = '  +0.28%  '

# Row 22 (WrappedliquidstakedEther2.0)
$ws.Range('D22').Value = '2.132.32'
$ws.Range('E22').Value = '  -0.29%  '

# Row 23 (Cosmos)
Set-TextValue $ws.Range('D23') '10.82'
$ws.Range('E23').Value = '  +0.24%  '

# Row 24 (Chainlink)
Set-TextValue $ws.Range('D24') '6.581'
$ws.Range('E24').Value = '  -0.04%  '

# Row 25 (Monero)
Set-TextValue $ws.Range('D25') '154.11'
$ws.Range('E25').Value = '  -0.10%  '

# Row 26 (Toncoin)
Set-TextValue $ws.Range('D26') '1.861'
$ws.Range('E26').Value = '  -2.60%  '

# Row 27 (LidoDAOToken)
Set-TextValue $ws.Range('D27') '2.182'
$ws.Range('E27').Value = '  +2.80%  '

# Row 28 (EthereumClassic)
Set-TextValue $ws.Range('D28') '18.30'
$ws.Range('E28').Value = '  -0.44%  '

# Row 29 (BitcoinCash)
Set-TextValue $ws.Range('D29') '114.74'
$ws.Range('E29').Value = '  -1.43%  '

# Row 30 (InternetComputer(DFINITY))
Set-TextValue $ws.Range('D30') '4.841'
$ws.Range('E30').Value = '  -1.56%  '

# Row 31 (Stellar)
Set-TextValue $ws.Range('D31') '0.08953'
$ws.Range('E31').Value = '  +0.14%  '

# Row 32 (HuobiToken)
Set-TextValue $ws.Range('D32') '3.183'
$ws.Range('E32').Value = '  +0.75%  '

# Row 33 (Filecoin)
$ws.Range('B33').Value = 'ARBITRUM'
$ws.Range('C33').Value = 'https://coinranking.com/coin/1Uo6s62Oc+arbitrum-arb'
Set-TextValue $ws.Range('D33') '1.233'
$ws.Range('E33').Value = '  +0.03%  '

# Row 34 (ARBITRUM)
$ws.Range('B34').Value = 'Filecoin'
$ws.Range('C34').Value = 'https://coinranking.com/coin/ymQub4fuB+filecoin-fil'
Set-TextValue $ws.Range('D34') '4.793'
$ws.Range('E34').Value = '  +3.55%  '

# Row 35 (ImmutableX)
Set-TextValue $ws.Range('D35') '0.7771'
$ws.Range('E35').Value = '  +2.35%  '

# Row 36 (VeChain)
$ws.Range('B36').Value = 'RenderToken'
$ws.Range('C36').Value = 'https://coinranking.com/coin/7C4Mh4xy1yDel+rendertoken-rndr'
Set-TextValue $ws.Range('D36') '2.616'
$ws.Range('E36').Value = '  +3.80%  '

# Row 37 (RenderToken)
$ws.Range('B37').Value = 'VeChain'
$ws.Range('C37').Value = 'https://coinranking.com/coin/FEbS54wxo4oIl+vechain-vet'
Set-TextValue $ws.Range('D37') '0.02088'
$ws.Range('E37').Value = '  +2.54%  '

# Row 38 (MXToken)
Set-TextValue $ws.Range('D38') '3.052'
$ws.Range('E38').Value = '  +2.48%  '

# Row 39 (TrustWalletToken)
$ws.Range('E39').Value = '  -0.27%  '

# Row 40 (TheSandbox)
Set-TextValue $ws.Range('D40') '0.5487'
$ws.Range('E40').Value = '  +1.08%  '

# Row 41 (Hedera)
Set-TextValue $ws.Range('D41') '0.05251'
$ws.Range('E41').Value = '  -0.08%  '

# Row 42 (FraxShare)
Set-TextValue $ws.Range('D42') '6.685'
$ws.Range('E42').Value = '  -3.68%  '

# Row 43 (Quant)
Set-TextValue $ws.Range('D43') '113.40'
$ws.Range('E43').Value = '  +4.36%  '

# Row 44 (Aptos)
Set-TextValue $ws.Range('D44') '8.477'
$ws.Range('E44').Value = '  +2.32%  '

# Row 45 (Algorand)
Set-TextValue $ws.Range('D45') '0.1504'
$ws.Range('E45').Value = '  -0.83%  '

# Row 46 (Decentraland)
Set-TextValue $ws.Range('D46') '0.4779'
$ws.Range('E46').Value = '  +0.07%  '

# Row 47 (EnergySwap)
Set-TextValue $ws.Range('D47') '10.42'
$ws.Range('E47').Value = '  -1.75%  '

# Row 48 (PaxDollar)
Set-TextValue $ws.Range('D48') '0.9998'
$ws.Range('E48').Value = '  -0.55%  '

# Row 49 (NEARProtocol)
Set-TextValue $ws.Range('D49') '1.611'
$ws.Range('E49').Value = '  -1.38%  '

# Row 50 (Aave)
Set-TextValue $ws.Range('D50') '66.60'
$ws.Range('E50').Value = '  -1.30%  '

# Row 51 (Cronos)
Set-TextValue $ws.Range('D51') '0.05995'
$ws.Range('E51').Value = '  -1.29%  '
